# Xbee sender BOM cleanup:
#  - Remove the extra "Header2/JTAG USER" connector row (old row 33)
#  - Remove the LogicAnalyzer (J6) and Mammond1593Y enclosure (Box1) rows (old rows 38-39)
#  - Remove the tantalum-capacitor (C2 C10 C12) and the three testpoint rows (old rows 42-45)
#  - Shift the remaining component rows up to close the gaps
#  - Add a subtotal (Price) formula and fix up the grand-total formulas

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Rows 33-41: re-write with the rows that remain, shifted upward ----

# Row 33 (was row 34): LEDT1.75 / SIP / D1
$ws.Range("A33").Value2 = "LEDT1.75"
$ws.Range("B33").Value2 = "SIP"
$ws.Range("D33").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("H33").Value2 = "D1"
$ws.Range("I33").Value2 = 1

# Row 34 (was row 35): LM3S811 / QFP / U1
$ws.Range("A34").Value2 = "LM3S811"
$ws.Range("B34").Value2 = "QFP"
$ws.Range("H34").Value2 = "U1"
$ws.Range("I34").Value2 = 1

# Row 35 (was row 36): LM2937ET-3.3 / TO220 / U2
$ws.Range("A35").Value2 = "LM2937ET-3.3"
$ws.Range("B35").Value2 = "TO220"
$ws.Range("H35").Value2 = "U2"
$ws.Range("I35").Value2 = 1
$ws.Range("K35").Value2 = 1.63

# Row 36 (was row 37): LM4041 / MISC / U4
$ws.Range("A36").Value2 = "LM4041"
$ws.Range("B36").Value2 = "MISC"
$ws.Range("H36").Value2 = "U4"
$ws.Range("I36").Value2 = 1
$ws.Range("K36").Value2 = 0.23

# Row 37 (was row 40): MC34119 / DIP / U5
$ws.Range("A37").Value2 = "MC34119"
$ws.Range("B37").Value2 = "DIP"
$ws.Range("H37").Value2 = "U5"
$ws.Range("I37").Value2 = 1
$ws.Range("J37").Value2 = "2.7 ma"
$ws.Range("K37").ClearContents()

# Row 38 (was row 41): Audio_Jack 3.5mm / CONN / CONN1
$ws.Range("A38").Value2 = "Audio_Jack 3.5mm"
$ws.Range("B38").Value2 = "CONN"
$ws.Range("H38").Value2 = "CONN1"
$ws.Range("I38").Value2 = 1

# Row 39 (was row 46): TLV5618 / DIP / U6
$ws.Range("A39").Value2 = "TLV5618"
$ws.Range("B39").Value2 = "DIP"
$ws.Range("H39").Value2 = "U6"
$ws.Range("I39").Value2 = 1

# Row 40 (was row 47): XTAL / DSC / Y1
$ws.Range("A40").Value2 = "XTAL"
$ws.Range("B40").Value2 = "DSC"
$ws.Range("H40").Value2 = "Y1"
$ws.Range("I40").Value2 = 1
$ws.Range("J40").ClearContents()

# Row 41 (was row 48): NHD-320240WG-BxTGK-VZ#-3VR-C LCD / Newhaven Display
$ws.Range("A41").Value2 = "NHD-320240WG-BxTGK-VZ#-3VR-C LCD"
$ws.Range("B41").ClearContents()
$ws.Range("D41").Value2 = "Newhaven Display"
$ws.Range("H41").ClearContents()
$ws.Range("I41").Value2 = 1
$ws.Range("K41").Value2 = 20

# ---- Row 42: clear stale leftover data, add Price subtotal formula ----
$ws.Range("A42").ClearContents()
$ws.Range("B42").ClearContents()
$ws.Range("H42").ClearContents()
$ws.Range("I42").ClearContents()
$ws.Range("K42").Formula = "=SUM(K2:K41)"

# ---- Rows 43-48: now empty, clear out remaining stale leftovers ----
$ws.Range("A43:K43").ClearContents()
$ws.Range("A44:K44").ClearContents()
$ws.Range("A45:K45").ClearContents()
$ws.Range("A46:K46").ClearContents()
$ws.Range("A47:K47").ClearContents()
$ws.Range("A48:K48").ClearContents()

# ---- Grand totals: ranges shrink to reflect the removed rows ----
$ws.Range("I50").Formula = "=SUM(I2:I45)"
$ws.Range("K50").Formula = "=SUM(K1:K45)"

# ---- View state: selection moves to K42 and the window scrolls down one row ----
$excel.ActiveWindow.ScrollRow = 23
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K42").Select() | Out-Null
